$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Begründung"
$ws.Range("B6").Value = "Zur Authentifizierung und als Backup falls der PC gestört ist"

$ws.Range("B6").Select()
